$wb = $excel.ActiveWorkbook

# LL_max_6 (sheet index 9)
$ws = $wb.Worksheets.Item("LL_max_6")
$ws.Range("O2").Value = 0.8909260930523819
$ws.Range("H3").Value = 8.082245025748676
$ws.Range("I3").Value = 4.690229336836681
$ws.Range("J3").Value = 0.1077319054061828
$ws.Range("K3").Value = 2.308438694630816
$ws.Range("L3").Value = 0.1077319054134128
$ws.Range("M3").Value = 2.308438694642946
$ws.Range("N3").Value = 0.9526279648020816
$ws.Range("Q3").Value = (-4.036469913634587*[Math]::Pow(10,-12))
$ws.Range("Q4").Value = (1.415091525375082*[Math]::Pow(10,-9))
$ws.Range("N5").Value = 0.9526279647663866
$ws.Range("Q5").Value = (1.888130047150138*[Math]::Pow(10,-9))
$ws.Range("S5").Value = 179.9999999981034
$ws.Range("N6").Value = 0.9526279647663866
$ws.Range("P6").Value = 0.952627964838214
$ws.Range("Q6").Value = (1.888130047150138*[Math]::Pow(10,-9))
$ws.Range("S6").Value = 179.9999999981034

# LL_max_10 (sheet index 10)
$ws = $wb.Worksheets.Item("LL_max_10")
$ws.Range("O2").Value = 0.8909260930523819
$ws.Range("H3").Value = 8.082245025748676
$ws.Range("I3").Value = 4.690229336836681
$ws.Range("J3").Value = 0.1077319054061828
$ws.Range("K3").Value = 2.308438694630816
$ws.Range("L3").Value = 0.1077319054134128
$ws.Range("M3").Value = 2.308438694642946
$ws.Range("N3").Value = 0.9526279648020816
$ws.Range("Q3").Value = (-4.036469913634587*[Math]::Pow(10,-12))
$ws.Range("Q4").Value = (1.415091525375082*[Math]::Pow(10,-9))
$ws.Range("N5").Value = 0.9526279647663866
$ws.Range("Q5").Value = (1.888130047150138*[Math]::Pow(10,-9))
$ws.Range("S5").Value = 179.9999999981034
$ws.Range("N6").Value = 0.9526279647663866
$ws.Range("P6").Value = 0.952627964838214
$ws.Range("Q6").Value = (1.888130047150138*[Math]::Pow(10,-9))
$ws.Range("S6").Value = 179.9999999981034

# LL_max_fault_6 (sheet index 11)
$ws = $wb.Worksheets.Item("LL_max_fault_6")
$ws.Range("H3").Value = 8.082245025748676
$ws.Range("I3").Value = 4.690229336836681
$ws.Range("J3").Value = 0.1077319054061828
$ws.Range("K3").Value = 2.308438694630816
$ws.Range("L3").Value = 0.1077319054134128
$ws.Range("M3").Value = 2.308438694642946
$ws.Range("O3").Value = 0.7109795688765317
$ws.Range("P4").Value = 0.9171656569531275
$ws.Range("R4").Value = -106.527958690445
$ws.Range("O5").Value = 0.7109795688758747
$ws.Range("Q5").Value = 17.92198046351095
$ws.Range("O6").Value = 0.7109795688758747
$ws.Range("P6").Value = 0.9171656569581229
$ws.Range("Q6").Value = 17.92198046351095

# LL_max_fault_10 (sheet index 12)
$ws = $wb.Worksheets.Item("LL_max_fault_10")
$ws.Range("H3").Value = 8.082245025748676
$ws.Range("I3").Value = 4.690229336836681
$ws.Range("J3").Value = 0.1077319054061828
$ws.Range("K3").Value = 2.308438694630816
$ws.Range("L3").Value = 0.1077319054134128
$ws.Range("M3").Value = 2.308438694642946
$ws.Range("O3").Value = 0.7109795688765317
$ws.Range("P4").Value = 0.9171656569531275
$ws.Range("R4").Value = -106.527958690445
$ws.Range("O5").Value = 0.7109795688758747
$ws.Range("Q5").Value = 17.92198046351095
$ws.Range("O6").Value = 0.7109795688758747
$ws.Range("P6").Value = 0.9171656569581229
$ws.Range("Q6").Value = 17.92198046351095

# LL_min_6 (sheet index 13)
$ws = $wb.Worksheets.Item("LL_min_6")
$ws.Range("N2").Value = 0.9477232070952269
$ws.Range("O2").Value = 0.7936674232739478
$ws.Range("H3").Value = 13.19853631846006
$ws.Range("I3").Value = 4.655449610515097
$ws.Range("J3").Value = 0.1153518594838952
$ws.Range("K3").Value = 2.416397540698012
$ws.Range("L3").Value = 0.1153518594714223
$ws.Range("M3").Value = 2.416397540675882
$ws.Range("N3").Value = 0.8660254037841224
$ws.Range("Q3").Value = (3.077617116946378*[Math]::Pow(10,-11))
$ws.Range("N4").Value = 0.8660254037597852
$ws.Range("Q4").Value = (2.75547625054083*[Math]::Pow(10,-9))
$ws.Range("N5").Value = 0.8660254037516727
$ws.Range("Q5").Value = (3.663704698727089*[Math]::Pow(10,-9))
$ws.Range("N6").Value = 0.8660254037516727
$ws.Range("Q6").Value = (3.663704698727089*[Math]::Pow(10,-9))

# LL_min_10 (sheet index 14)
$ws = $wb.Worksheets.Item("LL_min_10")
$ws.Range("N2").Value = 0.9477232070952269
$ws.Range("O2").Value = 0.7936674232739478
$ws.Range("H3").Value = 13.19853631846006
$ws.Range("I3").Value = 4.655449610515097
$ws.Range("J3").Value = 0.1153518594838952
$ws.Range("K3").Value = 2.416397540698012
$ws.Range("L3").Value = 0.1153518594714223
$ws.Range("M3").Value = 2.416397540675882
$ws.Range("N3").Value = 0.8660254037841224
$ws.Range("Q3").Value = (3.077617116946378*[Math]::Pow(10,-11))
$ws.Range("N4").Value = 0.8660254037597852
$ws.Range("Q4").Value = (2.75547625054083*[Math]::Pow(10,-9))
$ws.Range("N5").Value = 0.8660254037516727
$ws.Range("Q5").Value = (3.663704698727089*[Math]::Pow(10,-9))
$ws.Range("N6").Value = 0.8660254037516727
$ws.Range("Q6").Value = (3.663704698727089*[Math]::Pow(10,-9))

# LL_min_fault_6 (sheet index 15)
$ws = $wb.Worksheets.Item("LL_min_fault_6")
$ws.Range("P2").Value = 0.9648695716021195
$ws.Range("R2").Value = -92.12802105335015
$ws.Range("H3").Value = 13.19853631846006
$ws.Range("I3").Value = 4.655449610515097
$ws.Range("J3").Value = 0.1153518594838952
$ws.Range("K3").Value = 2.416397540698012
$ws.Range("L3").Value = 0.1153518594714223
$ws.Range("M3").Value = 2.416397540675882
$ws.Range("P3").Value = 0.8307352055776825
$ws.Range("Q4").Value = 17.57011779629442
$ws.Range("N5").Value = 1.005693895246945
$ws.Range("P5").Value = 0.8307352056047123
$ws.Range("R5").Value = -106.9885815710309
$ws.Range("N6").Value = 1.005693895246945
$ws.Range("P6").Value = 0.8307352056047123
$ws.Range("R6").Value = -106.9885815710309

# LL_min_fault_10 (sheet index 16)
$ws = $wb.Worksheets.Item("LL_min_fault_10")
$ws.Range("P2").Value = 0.9648695716021195
$ws.Range("R2").Value = -92.12802105335015
$ws.Range("H3").Value = 13.19853631846006
$ws.Range("I3").Value = 4.655449610515097
$ws.Range("J3").Value = 0.1153518594838952
$ws.Range("K3").Value = 2.416397540698012
$ws.Range("L3").Value = 0.1153518594714223
$ws.Range("M3").Value = 2.416397540675882
$ws.Range("P3").Value = 0.8307352055776825
$ws.Range("Q4").Value = 17.57011779629442
$ws.Range("N5").Value = 1.005693895246945
$ws.Range("P5").Value = 0.8307352056047123
$ws.Range("R5").Value = -106.9885815710309
$ws.Range("N6").Value = 1.005693895246945
$ws.Range("P6").Value = 0.8307352056047123
$ws.Range("R6").Value = -106.9885815710309

# LG_max_6 (sheet index 17)
$ws = $wb.Worksheets.Item("LG_max_6")
$ws.Range("N2").Value = 1.033161137394682
$ws.Range("O2").Value = 1.100000023845433
$ws.Range("P2").Value = 1.074419873344362
$ws.Range("Q2").Value = 29.60975627564915
$ws.Range("S2").Value = 146.7226792525623
$ws.Range("B3").Value = 3.05600557129472
$ws.Range("E3").Value = 35.28771278464006
$ws.Range("H3").Value = 8.082245025748684
$ws.Range("I3").Value = 4.690229336836486
$ws.Range("J3").Value = 0.1077319053183132
$ws.Range("K3").Value = 2.308438694642973
$ws.Range("L3").Value = 0.107731905413177
$ws.Range("M3").Value = 2.308438694642607
$ws.Range("N3").Value = 0.753327220157332
$ws.Range("O3").Value = 1.100000023849472
$ws.Range("P3").Value = 1.028046430022632
$ws.Range("Q3").Value = 25.7733279752319
$ws.Range("R3").Value = -89.99999999999754
$ws.Range("S3").Value = 131.2907152131199
$ws.Range("T3").Value = 3.05600557129472
$ws.Range("N4").Value = 0.7533272201681387
$ws.Range("O4").Value = 1.100000023849472
$ws.Range("P4").Value = 1.028046430006054
$ws.Range("Q4").Value = 25.77332797676842
$ws.Range("R4").Value = -89.9999999999978
$ws.Range("S4").Value = 131.2907152140013
$ws.Range("N5").Value = 0.7533272201717424
$ws.Range("O5").Value = 1.100000023849472
$ws.Range("P5").Value = 1.028046430000528
$ws.Range("Q5").Value = 25.77332797728053
$ws.Range("R5").Value = -89.99999999999793
$ws.Range("S5").Value = 131.2907152142951
$ws.Range("N6").Value = 0.7533272201717424
$ws.Range("O6").Value = 1.100000023849472
$ws.Range("P6").Value = 1.028046430000528
$ws.Range("Q6").Value = 25.77332797728053
$ws.Range("R6").Value = -89.99999999999793
$ws.Range("S6").Value = 131.2907152142951

# LG_max_10 (sheet index 18)
$ws = $wb.Worksheets.Item("LG_max_10")
$ws.Range("N2").Value = 1.033161137394682
$ws.Range("O2").Value = 1.100000023845433
$ws.Range("P2").Value = 1.074419873344362
$ws.Range("Q2").Value = 29.60975627564915
$ws.Range("S2").Value = 146.7226792525623
$ws.Range("B3").Value = 3.05600557129472
$ws.Range("E3").Value = 35.28771278464006
$ws.Range("H3").Value = 8.082245025748684
$ws.Range("I3").Value = 4.690229336836486
$ws.Range("J3").Value = 0.1077319053183132
$ws.Range("K3").Value = 2.308438694642973
$ws.Range("L3").Value = 0.107731905413177
$ws.Range("M3").Value = 2.308438694642607
$ws.Range("N3").Value = 0.753327220157332
$ws.Range("O3").Value = 1.100000023849472
$ws.Range("P3").Value = 1.028046430022632
$ws.Range("Q3").Value = 25.7733279752319
$ws.Range("R3").Value = -89.99999999999754
$ws.Range("S3").Value = 131.2907152131199
$ws.Range("T3").Value = 3.05600557129472
$ws.Range("N4").Value = 0.7533272201681387
$ws.Range("O4").Value = 1.100000023849472
$ws.Range("P4").Value = 1.028046430006054
$ws.Range("Q4").Value = 25.77332797676842
$ws.Range("R4").Value = -89.9999999999978
$ws.Range("S4").Value = 131.2907152140013
$ws.Range("N5").Value = 0.7533272201717424
$ws.Range("O5").Value = 1.100000023849472
$ws.Range("P5").Value = 1.028046430000528
$ws.Range("Q5").Value = 25.77332797728053
$ws.Range("R5").Value = -89.99999999999793
$ws.Range("S5").Value = 131.2907152142951
$ws.Range("N6").Value = 0.7533272201717424
$ws.Range("O6").Value = 1.100000023849472
$ws.Range("P6").Value = 1.028046430000528
$ws.Range("Q6").Value = 25.77332797728053
$ws.Range("R6").Value = -89.99999999999793
$ws.Range("S6").Value = 131.2907152142951

# LG_max_fault_6 (sheet index 19)
$ws = $wb.Worksheets.Item("LG_max_fault_6")
$ws.Range("N2").Value = 1.075354100191947
$ws.Range("O2").Value = 1.100000023843188
$ws.Range("P2").Value = 1.090947544103742
$ws.Range("Q2").Value = 29.81365834489262
$ws.Range("S2").Value = 148.7868404542143
$ws.Range("B3").Value = 1.131748224437117
$ws.Range("E3").Value = 13.06830284067301
$ws.Range("H3").Value = 8.082245025748684
$ws.Range("I3").Value = 4.690229336836486
$ws.Range("J3").Value = 0.1077319053183132
$ws.Range("K3").Value = 2.308438694642973
$ws.Range("L3").Value = 0.107731905413177
$ws.Range("M3").Value = 2.308438694642607
$ws.Range("N3").Value = 0.9720083825024672
$ws.Range("O3").Value = 1.100000023844735
$ws.Range("P3").Value = 1.064882470539624
$ws.Range("Q3").Value = 28.51404249107572
$ws.Range("R3").Value = -89.9999999999968
$ws.Range("S3").Value = 143.327655084583
$ws.Range("T3").Value = 1.131748224437117
$ws.Range("N4").Value = 0.9720083825065572
$ws.Range("O4").Value = 1.100000023844735
$ws.Range("P4").Value = 1.064882470534591
$ws.Range("Q4").Value = 28.51404249151392
$ws.Range("R4").Value = -89.99999999999692
$ws.Range("S4").Value = 143.3276550849505
$ws.Range("N5").Value = 0.9720083825079209
$ws.Range("O5").Value = 1.100000023844735
$ws.Range("P5").Value = 1.064882470532914
$ws.Range("Q5").Value = 28.51404249165996
$ws.Range("R5").Value = -89.99999999999696
$ws.Range("S5").Value = 143.3276550850731
$ws.Range("N6").Value = 0.9720083825079209
$ws.Range("O6").Value = 1.100000023844735
$ws.Range("P6").Value = 1.064882470532914
$ws.Range("Q6").Value = 28.51404249165996
$ws.Range("R6").Value = -89.99999999999696
$ws.Range("S6").Value = 143.3276550850731

# LG_max_fault_10 (sheet index 20)
$ws = $wb.Worksheets.Item("LG_max_fault_10")
$ws.Range("N2").Value = 1.075354100191947
$ws.Range("O2").Value = 1.100000023843188
$ws.Range("P2").Value = 1.090947544103742
$ws.Range("Q2").Value = 29.81365834489262
$ws.Range("S2").Value = 148.7868404542143
$ws.Range("B3").Value = 1.131748224437117
$ws.Range("E3").Value = 13.06830284067301
$ws.Range("H3").Value = 8.082245025748684
$ws.Range("I3").Value = 4.690229336836486
$ws.Range("J3").Value = 0.1077319053183132
$ws.Range("K3").Value = 2.308438694642973
$ws.Range("L3").Value = 0.107731905413177
$ws.Range("M3").Value = 2.308438694642607
$ws.Range("N3").Value = 0.9720083825024672
$ws.Range("O3").Value = 1.100000023844735
$ws.Range("P3").Value = 1.064882470539624
$ws.Range("Q3").Value = 28.51404249107572
$ws.Range("R3").Value = -89.9999999999968
$ws.Range("S3").Value = 143.327655084583
$ws.Range("T3").Value = 1.131748224437117
$ws.Range("N4").Value = 0.9720083825065572
$ws.Range("O4").Value = 1.100000023844735
$ws.Range("P4").Value = 1.064882470534591
$ws.Range("Q4").Value = 28.51404249151392
$ws.Range("R4").Value = -89.99999999999692
$ws.Range("S4").Value = 143.3276550849505
$ws.Range("N5").Value = 0.9720083825079209
$ws.Range("O5").Value = 1.100000023844735
$ws.Range("P5").Value = 1.064882470532914
$ws.Range("Q5").Value = 28.51404249165996
$ws.Range("R5").Value = -89.99999999999696
$ws.Range("S5").Value = 143.3276550850731
$ws.Range("N6").Value = 0.9720083825079209
$ws.Range("O6").Value = 1.100000023844735
$ws.Range("P6").Value = 1.064882470532914
$ws.Range("Q6").Value = 28.51404249165996
$ws.Range("R6").Value = -89.99999999999696
$ws.Range("S6").Value = 143.3276550850731

# LG_min_6 (sheet index 21)
$ws = $wb.Worksheets.Item("LG_min_6")
$ws.Range("N2").Value = 0.9503704895552483
$ws.Range("O2").Value = 0.9999999999968261
$ws.Range("P2").Value = 0.991342140678414
$ws.Range("Q2").Value = 28.96372540325381
$ws.Range("R2").Value = -89.99999999999639
$ws.Range("S2").Value = 147.0100636752413
$ws.Range("B3").Value = 2.106737123846299
$ws.Range("E3").Value = 24.32650491128878
$ws.Range("H3").Value = 13.19853631846007
$ws.Range("I3").Value = 4.65544961051488
$ws.Range("J3").Value = 0.115351859513187
$ws.Range("K3").Value = 2.416397540675832
$ws.Range("L3").Value = 0.1153518594710931
$ws.Range("M3").Value = 2.416397540675332
$ws.Range("N3").Value = 0.7697252297623911
$ws.Range("O3").Value = 0.9999999999993349
$ws.Range("P3").Value = 0.9971062261237972
$ws.Range("Q3").Value = 22.86845159258527
$ws.Range("R3").Value = -89.99999999999777
$ws.Range("S3").Value = 135.3393797755744
$ws.Range("T3").Value = 2.106737123846299
$ws.Range("N4").Value = 0.7697252297746899
$ws.Range("O4").Value = 0.9999999999993349
$ws.Range("P4").Value = 0.9971062261090393
$ws.Range("Q4").Value = 22.86845159415236
$ws.Range("R4").Value = -89.99999999999812
$ws.Range("S4").Value = 135.3393797766896
$ws.Range("N5").Value = 0.7697252297787915
$ws.Range("O5").Value = 0.9999999999993349
$ws.Range("P5").Value = 0.9971062261041209
$ws.Range("Q5").Value = 22.86845159467465
$ws.Range("R5").Value = -89.99999999999832
$ws.Range("S5").Value = 135.3393797770614
$ws.Range("N6").Value = 0.7697252297787915
$ws.Range("O6").Value = 0.9999999999993349
$ws.Range("P6").Value = 0.9971062261041209
$ws.Range("Q6").Value = 22.86845159467465
$ws.Range("R6").Value = -89.99999999999832
$ws.Range("S6").Value = 135.3393797770614

# LG_min_10 (sheet index 22)
$ws = $wb.Worksheets.Item("LG_min_10")
$ws.Range("N2").Value = 0.9503704895552483
$ws.Range("O2").Value = 0.9999999999968261
$ws.Range("P2").Value = 0.991342140678414
$ws.Range("Q2").Value = 28.96372540325381
$ws.Range("R2").Value = -89.99999999999639
$ws.Range("S2").Value = 147.0100636752413
$ws.Range("B3").Value = 2.106737123846299
$ws.Range("E3").Value = 24.32650491128878
$ws.Range("H3").Value = 13.19853631846007
$ws.Range("I3").Value = 4.65544961051488
$ws.Range("J3").Value = 0.115351859513187
$ws.Range("K3").Value = 2.416397540675832
$ws.Range("L3").Value = 0.1153518594710931
$ws.Range("M3").Value = 2.416397540675332
$ws.Range("N3").Value = 0.7697252297623911
$ws.Range("O3").Value = 0.9999999999993349
$ws.Range("P3").Value = 0.9971062261237972
$ws.Range("Q3").Value = 22.86845159258527
$ws.Range("R3").Value = -89.99999999999777
$ws.Range("S3").Value = 135.3393797755744
$ws.Range("T3").Value = 2.106737123846299
$ws.Range("N4").Value = 0.7697252297746899
$ws.Range("O4").Value = 0.9999999999993349
$ws.Range("P4").Value = 0.9971062261090393
$ws.Range("Q4").Value = 22.86845159415236
$ws.Range("R4").Value = -89.99999999999812
$ws.Range("S4").Value = 135.3393797766896
$ws.Range("N5").Value = 0.7697252297787915
$ws.Range("O5").Value = 0.9999999999993349
$ws.Range("P5").Value = 0.9971062261041209
$ws.Range("Q5").Value = 22.86845159467465
$ws.Range("R5").Value = -89.99999999999832
$ws.Range("S5").Value = 135.3393797770614
$ws.Range("N6").Value = 0.7697252297787915
$ws.Range("O6").Value = 0.9999999999993349
$ws.Range("P6").Value = 0.9971062261041209
$ws.Range("Q6").Value = 22.86845159467465
$ws.Range("R6").Value = -89.99999999999832
$ws.Range("S6").Value = 135.3393797770614

# LG_min_fault_6 (sheet index 23)
$ws = $wb.Worksheets.Item("LG_min_fault_6")
$ws.Range("N2").Value = 0.9775767415566581
$ws.Range("O2").Value = 0.9999999999987124
$ws.Range("P2").Value = 0.9937101329214633
$ws.Range("Q2").Value = 29.6830896164372
$ws.Range("S2").Value = 148.7235901980441
$ws.Range("B3").Value = 0.9232203391147255
$ws.Range("E3").Value = 10.66043022618449
$ws.Range("H3").Value = 13.19853631846007
$ws.Range("I3").Value = 4.65544961051488
$ws.Range("J3").Value = 0.115351859513187
$ws.Range("K3").Value = 2.416397540675832
$ws.Range("L3").Value = 0.1153518594710931
$ws.Range("M3").Value = 2.416397540675332
$ws.Range("N3").Value = 0.8935075228653694
$ws.Range("O3").Value = 0.9999999999998567
$ws.Range("P3").Value = 0.980016071475948
$ws.Range("Q3").Value = 27.9626568429412
$ws.Range("R3").Value = -89.99999999999692
$ws.Range("S3").Value = 143.6379170039678
$ws.Range("T3").Value = 0.9232203391147253
$ws.Range("N4").Value = 0.8935075228723632
$ws.Range("O4").Value = 0.9999999999998567
$ws.Range("P4").Value = 0.9800160714716953
$ws.Range("Q4").Value = 27.96265684345922
$ws.Range("R4").Value = -89.99999999999707
$ws.Range("S4").Value = 143.6379170045408
$ws.Range("N5").Value = 0.8935075228746953
$ws.Range("O5").Value = 0.9999999999998567
$ws.Range("P5").Value = 0.9800160714702781
$ws.Range("Q5").Value = 27.96265684363187
$ws.Range("R5").Value = -89.99999999999716
$ws.Range("S5").Value = 143.6379170047319
$ws.Range("N6").Value = 0.8935075228746953
$ws.Range("O6").Value = 0.9999999999998567
$ws.Range("P6").Value = 0.9800160714702781
$ws.Range("Q6").Value = 27.96265684363187
$ws.Range("R6").Value = -89.99999999999716
$ws.Range("S6").Value = 143.6379170047319

# LG_min_fault_10 (sheet index 24)
$ws = $wb.Worksheets.Item("LG_min_fault_10")
$ws.Range("N2").Value = 0.9775767415566581
$ws.Range("O2").Value = 0.9999999999987124
$ws.Range("P2").Value = 0.9937101329214633
$ws.Range("Q2").Value = 29.6830896164372
$ws.Range("S2").Value = 148.7235901980441
$ws.Range("B3").Value = 0.9232203391147255
$ws.Range("E3").Value = 10.66043022618449
$ws.Range("H3").Value = 13.19853631846007
$ws.Range("I3").Value = 4.65544961051488
$ws.Range("J3").Value = 0.115351859513187
$ws.Range("K3").Value = 2.416397540675832
$ws.Range("L3").Value = 0.1153518594710931
$ws.Range("M3").Value = 2.416397540675332
$ws.Range("N3").Value = 0.8935075228653694
$ws.Range("O3").Value = 0.9999999999998567
$ws.Range("P3").Value = 0.980016071475948
$ws.Range("Q3").Value = 27.9626568429412
$ws.Range("R3").Value = -89.99999999999692
$ws.Range("S3").Value = 143.6379170039678
$ws.Range("T3").Value = 0.9232203391147253
$ws.Range("N4").Value = 0.8935075228723632
$ws.Range("O4").Value = 0.9999999999998567
$ws.Range("P4").Value = 0.9800160714716953
$ws.Range("Q4").Value = 27.96265684345922
$ws.Range("R4").Value = -89.99999999999707
$ws.Range("S4").Value = 143.6379170045408
$ws.Range("N5").Value = 0.8935075228746953
$ws.Range("O5").Value = 0.9999999999998567
$ws.Range("P5").Value = 0.9800160714702781
$ws.Range("Q5").Value = 27.96265684363187
$ws.Range("R5").Value = -89.99999999999716
$ws.Range("S5").Value = 143.6379170047319
$ws.Range("N6").Value = 0.8935075228746953
$ws.Range("O6").Value = 0.9999999999998567
$ws.Range("P6").Value = 0.9800160714702781
$ws.Range("Q6").Value = 27.96265684363187
$ws.Range("R6").Value = -89.99999999999716
$ws.Range("S6").Value = 143.6379170047319

# LLG_max_6 (sheet index 25)
$ws = $wb.Worksheets.Item("LLG_max_6")
$ws.Range("N2").Value = 1.027429940985827
$ws.Range("O2").Value = 0.8909260930531102
$ws.Range("P2").Value = 1.050912456599734
$ws.Range("Q2").Value = 24.72496978084876
$ws.Range("R2").Value = -89.28600430312945
$ws.Range("S2").Value = 153.9738333667478
$ws.Range("C3").Value = 5.537121390632967
$ws.Range("D3").Value = 4.057466806025851
$ws.Range("F3").Value = 63.93717050835158
$ws.Range("G3").Value = 46.85159105373992
$ws.Range("H3").Value = 8.082245025748676
$ws.Range("I3").Value = 4.690229336836681
$ws.Range("J3").Value = 0.1077319054061828
$ws.Range("K3").Value = 2.308438694630816
$ws.Range("L3").Value = 0.1077319054134128
$ws.Range("M3").Value = 2.308438694642946
$ws.Range("N3").Value = 0.8886186586399299
$ws.Range("P3").Value = 0.8886186586419419
$ws.Range("Q3").Value = -5.563851380466366
$ws.Range("S3").Value = 174.4361486195176
$ws.Range("T3").Value = 1.90189366929511
$ws.Range("N4").Value = 0.8886186586090479
$ws.Range("P4").Value = 0.8886186586689288
$ws.Range("Q4").Value = -5.563851378677295
$ws.Range("S4").Value = 174.4361486186131
$ws.Range("N5").Value = 0.8886186585987529
$ws.Range("P5").Value = 0.888618658677924
$ws.Range("Q5").Value = -5.563851378080949
$ws.Range("S5").Value = 174.4361486183116
$ws.Range("N6").Value = 0.8886186585987529
$ws.Range("P6").Value = 0.888618658677924
$ws.Range("Q6").Value = -5.563851378080949
$ws.Range("S6").Value = 174.4361486183116

# LLG_max_10 (sheet index 26)
$ws = $wb.Worksheets.Item("LLG_max_10")
$ws.Range("N2").Value = 1.027429940985827
$ws.Range("O2").Value = 0.8909260930531102
$ws.Range("P2").Value = 1.050912456599734
$ws.Range("Q2").Value = 24.72496978084876
$ws.Range("R2").Value = -89.28600430312945
$ws.Range("S2").Value = 153.9738333667478
$ws.Range("C3").Value = 5.537121390632967
$ws.Range("D3").Value = 4.057466806025851
$ws.Range("F3").Value = 63.93717050835158
$ws.Range("G3").Value = 46.85159105373992
$ws.Range("H3").Value = 8.082245025748676
$ws.Range("I3").Value = 4.690229336836681
$ws.Range("J3").Value = 0.1077319054061828
$ws.Range("K3").Value = 2.308438694630816
$ws.Range("L3").Value = 0.1077319054134128
$ws.Range("M3").Value = 2.308438694642946
$ws.Range("N3").Value = 0.8886186586399299
$ws.Range("P3").Value = 0.8886186586419419
$ws.Range("Q3").Value = -5.563851380466366
$ws.Range("S3").Value = 174.4361486195176
$ws.Range("T3").Value = 1.90189366929511
$ws.Range("N4").Value = 0.8886186586090479
$ws.Range("P4").Value = 0.8886186586689288
$ws.Range("Q4").Value = -5.563851378677295
$ws.Range("S4").Value = 174.4361486186131
$ws.Range("N5").Value = 0.8886186585987529
$ws.Range("P5").Value = 0.888618658677924
$ws.Range("Q5").Value = -5.563851378080949
$ws.Range("S5").Value = 174.4361486183116
$ws.Range("N6").Value = 0.8886186585987529
$ws.Range("P6").Value = 0.888618658677924
$ws.Range("Q6").Value = -5.563851378080949
$ws.Range("S6").Value = 174.4361486183116

# LLG_max_fault_6 (sheet index 27)
$ws = $wb.Worksheets.Item("LLG_max_fault_6")
$ws.Range("N2").Value = 1.090139876245251
$ws.Range("O2").Value = 1.052967690033645
$ws.Range("P2").Value = 1.073842909852251
$ws.Range("Q2").Value = 28.44553889965567
$ws.Range("R2").Value = -91.44286082238069
$ws.Range("S2").Value = 150.2182633220348
$ws.Range("C3").Value = 1.42282178157132
$ws.Range("D3").Value = 1.203761473393288
$ws.Range("F3").Value = 16.42933077198129
$ws.Range("G3").Value = 13.8998402140743
$ws.Range("H3").Value = 8.082245025748676
$ws.Range("I3").Value = 4.690229336836681
$ws.Range("J3").Value = 0.1077319054061828
$ws.Range("K3").Value = 2.308438694630816
$ws.Range("L3").Value = 0.1077319054134128
$ws.Range("M3").Value = 2.308438694642946
$ws.Range("N3").Value = 1.066180738057245
$ws.Range("O3").Value = 0.8723446855805088
$ws.Range("P3").Value = 0.9617364120304487
$ws.Range("Q3").Value = 21.48012101510393
$ws.Range("R3").Value = -100.0510163396169
$ws.Range("S3").Value = 150.844226026951
$ws.Range("T3").Value = 0.9268468105291111
$ws.Range("N4").Value = 1.066180738050239
$ws.Range("O4").Value = 0.8723446855802062
$ws.Range("P4").Value = 0.9617364120371624
$ws.Range("Q4").Value = 21.48012101553641
$ws.Range("R4").Value = -100.051016338402
$ws.Range("S4").Value = 150.8442260273101
$ws.Range("N5").Value = 1.066180738047904
$ws.Range("O5").Value = 0.8723446855801055
$ws.Range("P5").Value = 0.9617364120394001
$ws.Range("Q5").Value = 21.48012101568058
$ws.Range("R5").Value = -100.051016337997
$ws.Range("S5").Value = 150.8442260274298
$ws.Range("N6").Value = 1.066180738047904
$ws.Range("O6").Value = 0.8723446855801055
$ws.Range("P6").Value = 0.9617364120394001
$ws.Range("Q6").Value = 21.48012101568058
$ws.Range("R6").Value = -100.051016337997
$ws.Range("S6").Value = 150.8442260274298

# LLG_max_fault_10 (sheet index 28)
$ws = $wb.Worksheets.Item("LLG_max_fault_10")
$ws.Range("N2").Value = 1.090139876245251
$ws.Range("O2").Value = 1.052967690033645
$ws.Range("P2").Value = 1.073842909852251
$ws.Range("Q2").Value = 28.44553889965567
$ws.Range("R2").Value = -91.44286082238069
$ws.Range("S2").Value = 150.2182633220348
$ws.Range("C3").Value = 1.42282178157132
$ws.Range("D3").Value = 1.203761473393288
$ws.Range("F3").Value = 16.42933077198129
$ws.Range("G3").Value = 13.8998402140743
$ws.Range("H3").Value = 8.082245025748676
$ws.Range("I3").Value = 4.690229336836681
$ws.Range("J3").Value = 0.1077319054061828
$ws.Range("K3").Value = 2.308438694630816
$ws.Range("L3").Value = 0.1077319054134128
$ws.Range("M3").Value = 2.308438694642946
$ws.Range("N3").Value = 1.066180738057245
$ws.Range("O3").Value = 0.8723446855805088
$ws.Range("P3").Value = 0.9617364120304487
$ws.Range("Q3").Value = 21.48012101510393
$ws.Range("R3").Value = -100.0510163396169
$ws.Range("S3").Value = 150.844226026951
$ws.Range("T3").Value = 0.9268468105291111
$ws.Range("N4").Value = 1.066180738050239
$ws.Range("O4").Value = 0.8723446855802062
$ws.Range("P4").Value = 0.9617364120371624
$ws.Range("Q4").Value = 21.48012101553641
$ws.Range("R4").Value = -100.051016338402
$ws.Range("S4").Value = 150.8442260273101
$ws.Range("N5").Value = 1.066180738047904
$ws.Range("O5").Value = 0.8723446855801055
$ws.Range("P5").Value = 0.9617364120394001
$ws.Range("Q5").Value = 21.48012101568058
$ws.Range("R5").Value = -100.051016337997
$ws.Range("S5").Value = 150.8442260274298
$ws.Range("N6").Value = 1.066180738047904
$ws.Range("O6").Value = 0.8723446855801055
$ws.Range("P6").Value = 0.9617364120394001
$ws.Range("Q6").Value = 21.48012101568058
$ws.Range("R6").Value = -100.051016337997
$ws.Range("S6").Value = 150.8442260274298

# LLG_min_6 (sheet index 29)
$ws = $wb.Worksheets.Item("LLG_min_6")
$ws.Range("N2").Value = 0.9356305267033488
$ws.Range("O2").Value = 0.7936674232730628
$ws.Range("P2").Value = 0.9562742727128652
$ws.Range("Q2").Value = 24.21813576263254
$ws.Range("R2").Value = -89.2249262389752
$ws.Range("S2").Value = 154.6259603020324
$ws.Range("C3").Value = 4.675587159324556
$ws.Range("D3").Value = 3.609362078175536
$ws.Range("F3").Value = 53.98903010111181
$ws.Range("G3").Value = 41.67732334874946
$ws.Range("H3").Value = 13.19853631846006
$ws.Range("I3").Value = 4.655449610515097
$ws.Range("J3").Value = 0.1153518594838952
$ws.Range("K3").Value = 2.416397540698012
$ws.Range("L3").Value = 0.1153518594714223
$ws.Range("M3").Value = 2.416397540675882
$ws.Range("N3").Value = 0.836174233327264
$ws.Range("P3").Value = 0.836174233327136
$ws.Range("Q3").Value = -4.432341270692252
$ws.Range("S3").Value = 175.5676587292177
$ws.Range("T3").Value = 1.194917017012807
$ws.Range("N4").Value = 0.8361742332987967
$ws.Range("P4").Value = 0.8361742333535674
$ws.Range("Q4").Value = -4.43234126753898
$ws.Range("S4").Value = 175.5676587270017
$ws.Range("N5").Value = 0.8361742332893063
$ws.Range("P5").Value = 0.8361742333623772
$ws.Range("Q5").Value = -4.432341266487897
$ws.Range("S5").Value = 175.567658726263
$ws.Range("N6").Value = 0.8361742332893063
$ws.Range("P6").Value = 0.8361742333623772
$ws.Range("Q6").Value = -4.432341266487897
$ws.Range("S6").Value = 175.567658726263

# LLG_min_10 (sheet index 30)
$ws = $wb.Worksheets.Item("LLG_min_10")
$ws.Range("N2").Value = 0.9356305267033488
$ws.Range("O2").Value = 0.7936674232730628
$ws.Range("P2").Value = 0.9562742727128652
$ws.Range("Q2").Value = 24.21813576263254
$ws.Range("R2").Value = -89.2249262389752
$ws.Range("S2").Value = 154.6259603020324
$ws.Range("C3").Value = 4.675587159324556
$ws.Range("D3").Value = 3.609362078175536
$ws.Range("F3").Value = 53.98903010111181
$ws.Range("G3").Value = 41.67732334874946
$ws.Range("H3").Value = 13.19853631846006
$ws.Range("I3").Value = 4.655449610515097
$ws.Range("J3").Value = 0.1153518594838952
$ws.Range("K3").Value = 2.416397540698012
$ws.Range("L3").Value = 0.1153518594714223
$ws.Range("M3").Value = 2.416397540675882
$ws.Range("N3").Value = 0.836174233327264
$ws.Range("P3").Value = 0.836174233327136
$ws.Range("Q3").Value = -4.432341270692252
$ws.Range("S3").Value = 175.5676587292177
$ws.Range("T3").Value = 1.194917017012807
$ws.Range("N4").Value = 0.8361742332987967
$ws.Range("P4").Value = 0.8361742333535674
$ws.Range("Q4").Value = -4.43234126753898
$ws.Range("S4").Value = 175.5676587270017
$ws.Range("N5").Value = 0.8361742332893063
$ws.Range("P5").Value = 0.8361742333623772
$ws.Range("Q5").Value = -4.432341266487897
$ws.Range("S5").Value = 175.567658726263
$ws.Range("N6").Value = 0.8361742332893063
$ws.Range("P6").Value = 0.8361742333623772
$ws.Range("Q6").Value = -4.432341266487897
$ws.Range("S6").Value = 175.567658726263

# LLG_min_fault_6 (sheet index 31)
$ws = $wb.Worksheets.Item("LLG_min_fault_6")
$ws.Range("N2").Value = 0.991784009114929
$ws.Range("O2").Value = 0.9517844802523826
$ws.Range("P2").Value = 0.9749657938312043
$ws.Range("Q2").Value = 28.21114579926146
$ws.Range("R2").Value = -91.61401138281248
$ws.Range("S2").Value = 150.3329223236375
$ws.Range("C3").Value = 1.286353511802591
$ws.Range("D3").Value = 1.027329872979415
$ws.Range("F3").Value = 14.85353092624493
$ws.Range("G3").Value = 11.86258357422419
$ws.Range("H3").Value = 13.19853631846006
$ws.Range("I3").Value = 4.655449610515097
$ws.Range("J3").Value = 0.1153518594838952
$ws.Range("K3").Value = 2.416397540698012
$ws.Range("L3").Value = 0.1153518594714223
$ws.Range("M3").Value = 2.416397540675882
$ws.Range("N3").Value = 0.9776016549893407
$ws.Range("O3").Value = 0.7848521014138666
$ws.Range("P3").Value = 0.8777903335451241
$ws.Range("Q3").Value = 21.06843963439565
$ws.Range("R3").Value = -100.4045966786942
$ws.Range("S3").Value = 151.3758615068555
$ws.Range("T3").Value = 0.7016277373189992
$ws.Range("N4").Value = 0.9776016549826257
$ws.Range("O4").Value = 0.7848521014214425
$ws.Range("P4").Value = 0.8777903335563334
$ws.Range("Q4").Value = 21.06843963517454
$ws.Range("R4").Value = -100.4045966769021
$ws.Range("S4").Value = 151.3758615071137
$ws.Range("N5").Value = 0.9776016549803869
$ws.Range("O5").Value = 0.7848521014239676
$ws.Range("P5").Value = 0.8777903335600694
$ws.Range("Q5").Value = 21.06843963543418
$ws.Range("R5").Value = -100.4045966763046
$ws.Range("S5").Value = 151.3758615071997
$ws.Range("N6").Value = 0.9776016549803869
$ws.Range("O6").Value = 0.7848521014239676
$ws.Range("P6").Value = 0.8777903335600694
$ws.Range("Q6").Value = 21.06843963543418
$ws.Range("R6").Value = -100.4045966763046
$ws.Range("S6").Value = 151.3758615071997

# LLG_min_fault_10 (sheet index 32)
$ws = $wb.Worksheets.Item("LLG_min_fault_10")
$ws.Range("N2").Value = 0.991784009114929
$ws.Range("O2").Value = 0.9517844802523826
$ws.Range("P2").Value = 0.9749657938312043
$ws.Range("Q2").Value = 28.21114579926146
$ws.Range("R2").Value = -91.61401138281248
$ws.Range("S2").Value = 150.3329223236375
$ws.Range("C3").Value = 1.286353511802591
$ws.Range("D3").Value = 1.027329872979415
$ws.Range("F3").Value = 14.85353092624493
$ws.Range("G3").Value = 11.86258357422419
$ws.Range("H3").Value = 13.19853631846006
$ws.Range("I3").Value = 4.655449610515097
$ws.Range("J3").Value = 0.1153518594838952
$ws.Range("K3").Value = 2.416397540698012
$ws.Range("L3").Value = 0.1153518594714223
$ws.Range("M3").Value = 2.416397540675882
$ws.Range("N3").Value = 0.9776016549893407
$ws.Range("O3").Value = 0.7848521014138666
$ws.Range("P3").Value = 0.8777903335451241
$ws.Range("Q3").Value = 21.06843963439565
$ws.Range("R3").Value = -100.4045966786942
$ws.Range("S3").Value = 151.3758615068555
$ws.Range("T3").Value = 0.7016277373189992
$ws.Range("N4").Value = 0.9776016549826257
$ws.Range("O4").Value = 0.7848521014214425
$ws.Range("P4").Value = 0.8777903335563334
$ws.Range("Q4").Value = 21.06843963517454
$ws.Range("R4").Value = -100.4045966769021
$ws.Range("S4").Value = 151.3758615071137
$ws.Range("N5").Value = 0.9776016549803869
$ws.Range("O5").Value = 0.7848521014239676
$ws.Range("P5").Value = 0.8777903335600694
$ws.Range("Q5").Value = 21.06843963543418
$ws.Range("R5").Value = -100.4045966763046
$ws.Range("S5").Value = 151.3758615071997
$ws.Range("N6").Value = 0.9776016549803869
$ws.Range("O6").Value = 0.7848521014239676
$ws.Range("P6").Value = 0.8777903335600694
$ws.Range("Q6").Value = 21.06843963543418
$ws.Range("R6").Value = -100.4045966763046
$ws.Range("S6").Value = 151.3758615071997

